$wb = $excel.ActiveWorkbook

$overviewSheet = $wb.Worksheets.Item("Overview")
$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

# Update the shared "Ready for handoff" -> "Handback transform failed" status
# (affects every sheet that shows this status for the 83c8cb06... file)
$overviewSheet.Cells.Item(3, 2).Value = "Handback transform failed"
$overviewSheet.Cells.Item(3, 3).Value = "Handback transform failed"
$zhSheet.Cells.Item(3, 3).Value = "Handback transform failed"
$deSheet.Cells.Item(3, 3).Value = "Handback transform failed"

# Add Error Detail (column L) messages for the failed handback transform
$zhSheet.Cells.Item(3, 12).Value = "Handback file name: qsda3lxl.1sa is different with handoff file name: 83c8cb06-4ef5-4895-82bb-5658c602c90d.7d8598339d2d89259d0127a768d04eca10ca6e3b.zh-cn."
$deSheet.Cells.Item(3, 12).Value = "Handback file name: qsda3lxl.1sa is different with handoff file name: 83c8cb06-4ef5-4895-82bb-5658c602c90d.7d8598339d2d89259d0127a768d04eca10ca6e3b.de-de."
